$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13 (Wish List): text stays the same, F13 becomes numeric 0 instead of "yet to test." ---
$ws.Range("F13").Value = 0

# --- Row 15 (was "(TS_008) Check Out") -> becomes "(TS_009) Check Out"; F15 becomes numeric 0 ---
$ws.Range("B15").Value = "(TS_009)`n Check Out"
$ws.Range("F15").Value = 0

# --- New row 16: "(TS_010) Header Desktop" ---
# Copy per-cell formatting from row 15 so the new row matches the look of the
# existing test-scenario rows (centered bold header cell, wrapped description, etc.)
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("F15").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D16").Value = "Validate the functionality of home page > Header > Desktop."

# --- Row 14 (was "(TS_007) Shopping Cart") -> becomes "(TS_008) Shopping Cart"; F14 becomes numeric 0 ---
$ws.Range("B14").Value = "(TS_008)`n Shopping Cart"
$ws.Range("F14").Value = 0

$ws.Range("B16").Value = "(TS_010)`n Header Desktop"
$ws.Range("C16").Value = "FRS"
$ws.Range("F16").Value = 16
$ws.Rows(16).RowHeight = 30

# --- View: scroll down a bit and move the active selection to B17 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B17").Select()
